# BAJAJ-PL MIS Base Page... (Performance_IDFC_HL.xlsx)
#
# Refreshes the BKT/STATE performance roll-up on Sheet1:
#   - rows 2-6 get new bucket labels, states and metrics
#   - two new rows (8 and 9 -> sheet rows 7-9 below the header) are added
#     for additional RAJASTHAN buckets
#   - the used range grows from A1:AE6 to A1:AE9 automatically as the
#     new cells are written

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order: BKT, STATE, TOTAL_POS, COUNT, FLOW_CASES, SB_CASES, RB_CASES,
# NM_CASES, PP_CASES, FC_CASES, SC_CASES, LOAN_CANCELED_CASES, FLOW_POS,
# SB_POS, RB_POS, NM_POS, PP_POS, FC_POS, SC_POS, LOAN_CANCELED_POS,
# FLOW_POS%, SB_POS%, RB_POS%, FC_POS%, SC_POS%, NM_POS%, PP_POS%,
# LOAN_CANCELED%, TOTAL PAID, POS_RES%, Additional_Performance
# i.e. columns A .. AE, one array per worksheet row (rows 2-9).
$rows = @(
    @("BKT0", "DELHI NCR", 817120112.4299999, 228, 25, 202, 0, 0, 0, 0, 1, 0, 61842820.34, 752412783.87, 0, 0, 0, 0, 2864508.22, 0, 7.57, 92.08, 0, 0, 0.35, 0, 0, 0, 8915728, 92.43000000000001, 0.35),
    @("BKT1", "DELHI NCR", 288503013.31, 199, 27, 125, 0, 44, 0, 0, 3, 0, 39913237.29, 194592171.52, 0, 50315521.08, 0, 0, 3682083.42, 0, 13.83, 67.45, 0, 0, 1.28, 17.44, 0, 0, 5234507, 86.17, 18.72),
    @("BKT1", "RAJASTHAN", 9911622.050000001, 15, 15, 0, 0, 0, 0, 0, 0, 0, 9911622.050000001, 0, 0, 0, 0, 0, 0, 0, 100, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @("BKT2", "RAJASTHAN", 1921240, 2, 2, 0, 0, 0, 0, 0, 0, 0, 1921240, 0, 0, 0, 0, 0, 0, 0, 100, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @("BKT3", "RAJASTHAN", 1765526.64, 4, 4, 0, 0, 0, 0, 0, 0, 0, 1765526.64, 0, 0, 0, 0, 0, 0, 0, 100, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @("BKT4", "RAJASTHAN", 489693, 2, 2, 0, 0, 0, 0, 0, 0, 0, 489693, 0, 0, 0, 0, 0, 0, 0, 100, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @("BKT5", "RAJASTHAN", 2298733, 2, 2, 0, 0, 0, 0, 0, 0, 0, 2298733, 0, 0, 0, 0, 0, 0, 0, 100, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @("BKT7", "RAJASTHAN", 1180938.14, 1, 1, 0, 0, 0, 0, 0, 0, 0, 1180938.14, 0, 0, 0, 0, 0, 0, 0, 100, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowData = $rows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
